$d = $word.ActiveDocument

function Replace-Text($searchText, $replaceText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $r.Text = $replaceText
    }
}

Replace-Text "each player has two options shown next to them, split or steal" "chaque joueur a deux options affichées à côté d'eux, diviser ou voler"

Replace-Text "Since two players have two choices each there are four outcomes in total" "Étant donné que deux joueurs ont deux choix chacun il y a quatre résultats au total"

Replace-Text "and they are all shown in the table in each section" "et ils sont tous affichés sur le tableau dans chaque section"

Replace-Text "The red number is the number of points won by the red player" "Le nombre en rouge est le nombre de points gagnés par le joueur rouge"

Replace-Text "and the blue number is the number of points won by the blue player" "et le nombre en bleu est le nombre de points gagnés par le joueur bleu"

# "for example" occurs three times in this subtitle document (it is a stock
# filler phrase used in several places), but only the single occurrence that
# immediately follows the "blue number" caption (and precedes "if both
# players choose to split") should be translated here. Locate that specific
# paragraph instead of doing a blanket Content.Find, so the other two
# occurrences of "for example" are left untouched.
$blueIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "et le nombre en bleu est le nombre de points gagnés par le joueur bleu") {
        $blueIdx = $i
    }
}

$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($blueIdx -gt 0 -and $i -gt $blueIdx) {
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq "for example") {
            $pr = $p.Range
            $found = $pr.Find.Execute("for example", $true, $false, $false, $false, $false, $true, 1, $false)
            if ($found) {
                $pr.Text = "par exemple,"
            }
            break
        }
    }
}

Replace-Text "if both players choose to split" "si les deux joueurs choisissent de diviser"

Replace-Text "we would end up with the top left outcome" "nous nous retrouverions avec le résultat en haut à gauche"
